# Apply "Update on 2018-07-03, 支出生活费400" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 currently has the same (blank-template) formatting as rows 45-50;
# pull in the real formatting used by the data rows above (row 43) first.
$ws.Range("D43:G43").Copy()
$ws.Range("D44:G44").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 44: new expense entry for 生活费 (7/1-7/10), 400 RMB, dated 2018-07-01.
$ws.Range("B44").Value = 42
$ws.Range("C44").Value = "支出"
$ws.Range("D44").Value = 400
$ws.Range("E44").Value = [DateTime]"2018-07-01"
$ws.Range("F44").Value = "生活费"
$ws.Range("G44").Value = "生活费(7/1-7/10)"

# Match the view state captured in the saved workbook (scrolled/selected cell).
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("H50").Select()
